$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-added "market aroma" translation rows (48-78) on the
# "Translations - Common" sheet. Column A is always the "cs" locale,
# column B is the translation key, column C is the Czech value.
# (Row 74 writes C before B to reproduce the shared-string insertion
# order of the original edit.)

$ws.Cells.Item(48,1).Value = "cs"
$ws.Cells.Item(48,2).Value = "market.index.title"
$ws.Cells.Item(48,3).Value = "Tržiště"

$ws.Cells.Item(49,1).Value = "cs"
$ws.Cells.Item(49,2).Value = "market.liquid.menu"
$ws.Cells.Item(49,3).Value = "Liquidy"

$ws.Cells.Item(50,1).Value = "cs"
$ws.Cells.Item(50,2).Value = "market.aroma.menu"
$ws.Cells.Item(50,3).Value = "Aromata"

$ws.Cells.Item(51,1).Value = "cs"
$ws.Cells.Item(51,2).Value = "market.build.menu"
$ws.Cells.Item(51,3).Value = "Buildy"

$ws.Cells.Item(52,1).Value = "cs"
$ws.Cells.Item(52,2).Value = "market.cotton.menu"
$ws.Cells.Item(52,3).Value = "Vaty"

$ws.Cells.Item(53,1).Value = "cs"
$ws.Cells.Item(53,2).Value = "market.hardware.menu"
$ws.Cells.Item(53,3).Value = "Hardware"

$ws.Cells.Item(54,1).Value = "cs"
$ws.Cells.Item(54,2).Value = "market.atomizer.menu"
$ws.Cells.Item(54,3).Value = "Atomizéry"

$ws.Cells.Item(55,1).Value = "cs"
$ws.Cells.Item(55,2).Value = "market.mod.menu"
$ws.Cells.Item(55,3).Value = "Mody"

$ws.Cells.Item(56,1).Value = "cs"
$ws.Cells.Item(56,2).Value = "market.cell.menu"
$ws.Cells.Item(56,3).Value = "Články"

$ws.Cells.Item(57,1).Value = "cs"
$ws.Cells.Item(57,2).Value = "market.other.menu"
$ws.Cells.Item(57,3).Value = "Ostatní"

$ws.Cells.Item(58,1).Value = "cs"
$ws.Cells.Item(58,2).Value = "market.voucher.menu"
$ws.Cells.Item(58,3).Value = "Poukázky"

$ws.Cells.Item(59,1).Value = "cs"
$ws.Cells.Item(59,2).Value = "market.vendor.menu"
$ws.Cells.Item(59,3).Value = "Výrobci"

$ws.Cells.Item(60,1).Value = "cs"
$ws.Cells.Item(60,2).Value = "market.certificate.menu"
$ws.Cells.Item(60,3).Value = "Certifikáty"

$ws.Cells.Item(61,1).Value = "cs"
$ws.Cells.Item(61,2).Value = "market.license.menu"
$ws.Cells.Item(61,3).Value = "Licence"

$ws.Cells.Item(62,1).Value = "cs"
$ws.Cells.Item(62,2).Value = "market.wire.menu"
$ws.Cells.Item(62,3).Value = "Odporové dráty"

$ws.Cells.Item(63,1).Value = "cs"
$ws.Cells.Item(63,2).Value = "inventory.index.title"
$ws.Cells.Item(63,3).Value = "Inventář"

$ws.Cells.Item(64,1).Value = "cs"
$ws.Cells.Item(64,2).Value = "market.aroma.index.title"
$ws.Cells.Item(64,3).Value = "Aromata"

$ws.Cells.Item(65,1).Value = "cs"
$ws.Cells.Item(65,2).Value = "common.infinite.loading"
$ws.Cells.Item(65,3).Value = "Přemýšlím"

$ws.Cells.Item(66,1).Value = "cs"
$ws.Cells.Item(66,2).Value = "common.infinite.no-more"
$ws.Cells.Item(66,3).Value = "Konec"

$ws.Cells.Item(67,1).Value = "cs"
$ws.Cells.Item(67,2).Value = "market.aroma.aroma.title"
$ws.Cells.Item(67,3).Value = "Náhled aromatu"

$ws.Cells.Item(68,1).Value = "cs"
$ws.Cells.Item(68,2).Value = "aroma.info.tab"
$ws.Cells.Item(68,3).Value = "Aroma"

$ws.Cells.Item(69,1).Value = "cs"
$ws.Cells.Item(69,2).Value = "aroma.more.tab"
$ws.Cells.Item(69,3).Value = "Více"

$ws.Cells.Item(70,1).Value = "cs"
$ws.Cells.Item(70,2).Value = "market.aroma.view.steep"
$ws.Cells.Item(70,3).Value = "Doba zrání"

$ws.Cells.Item(71,1).Value = "cs"
$ws.Cells.Item(71,2).Value = "market.aroma.view.tastes"
$ws.Cells.Item(71,3).Value = "Příchutě"

$ws.Cells.Item(72,1).Value = "cs"
$ws.Cells.Item(72,2).Value = "common.taste.empty"
$ws.Cells.Item(72,3).Value = "Bez uvedených příchutí"

$ws.Cells.Item(73,1).Value = "cs"
$ws.Cells.Item(73,2).Value = "market.aroma.view.name"
$ws.Cells.Item(73,3).Value = "Název"

$ws.Cells.Item(74,1).Value = "cs"
$ws.Cells.Item(74,3).Value = "Poměr VG/PG"
$ws.Cells.Item(74,2).Value = "market.aroma.view.vgpg"

$ws.Cells.Item(75,1).Value = "cs"
$ws.Cells.Item(75,2).Value = "common.pgvg.pg.tooltip"
$ws.Cells.Item(75,3).Value = "Poměr VG/PG"

$ws.Cells.Item(76,1).Value = "cs"
$ws.Cells.Item(76,2).Value = "market.aroma.view.content"
$ws.Cells.Item(76,3).Value = "Obsah aromatu"

$ws.Cells.Item(77,1).Value = "cs"
$ws.Cells.Item(77,2).Value = "common.aroma.content.tooltip"
$ws.Cells.Item(77,3).Value = "Obsah aromatu v lahvičce."

$ws.Cells.Item(78,1).Value = "cs"
$ws.Cells.Item(78,2).Value = "common.aroma.volume.tooltip"
$ws.Cells.Item(78,3).Value = "Objem lahvičky aromatu (v případě SnV), případně doporučený objem míchání pro čisté aroma."

# The sheet used to have one extra trailing blank row (1603); remove it so
# the used range ends at row 1602 again.
$ws.Rows.Item(1603).Delete()

# Move the cursor / selection to match the saved view state.
$ws.Range("B64").Select()
